# Updates the "cryptos" price/volume table for the Thu Mar 30 10:08:45 UTC 2023
# GitHub Actions refresh: new Price (col D) / Volume(1h) (col E) readings, plus a
# rank swap between Aptos and TrustWalletToken (rows 40/41, cols B/C/D/E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text assignment - safe for values Excel will not mis-parse as a number
# (links, coin names, percentage strings with padding, multi-dot price strings).
function Set-PlainText($addr, $value) {
    $ws.Range($addr).Value = $value
}

# The sheet stores "Price" (col D) as plain text, even for values that look like
# plain decimals (e.g. "1.001"). A bare Range.Value assignment would let Excel
# auto-convert those to numbers, so force text the same way Excel's own
# apostrophe-prefix entry does: flip to a text format, assign, then restore the
# original "Normal"/General formatting so no stray style is left behind.
function Set-TextForced($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-PlainText "D2" '28.628.45'
Set-PlainText "E2" '  +0.92%  '

Set-PlainText "D3" '1.803.60'
Set-PlainText "E3" '  -0.68%  '

Set-PlainText "E4" '  +0.07%  '

Set-TextForced "D5" '317.03'
Set-PlainText "E5" '  -0.37%  '

Set-TextForced "D6" '1.001'
Set-PlainText "E6" '  +0.05%  '

Set-TextForced "D7" '0.5314'
Set-PlainText "E7" '  -6.53%  '

Set-TextForced "D8" '0.3760'
Set-PlainText "E8" '  -3.14%  '

Set-TextForced "D9" '0.07505'

Set-TextForced "D10" '42.45'
Set-PlainText "E10" '  -1.66%  '

Set-PlainText "E11" '  -2.35%  '

Set-PlainText "E12" '  +0.09%  '

Set-TextForced "D13" '20.70'
Set-PlainText "E13" '  -2.53%  '

Set-TextForced "D14" '6.149'
Set-PlainText "E14" '  -1.73%  '

Set-TextForced "D15" '7.363'
Set-PlainText "E15" '  +0.77%  '

Set-PlainText "D16" '1.802.11'
Set-PlainText "E16" '  -0.56%  '

Set-TextForced "D17" '90.15'
Set-PlainText "E17" '  -2.11%  '

Set-PlainText "E18" '  -1.30%  '

Set-TextForced "D19" '0.06459'
Set-PlainText "E19" '  -0.35%  '

Set-TextForced "D20" '1.001'
Set-PlainText "E20" '  +0.06%  '

Set-TextForced "D21" '17.26'
Set-PlainText "E21" '  -0.47%  '

Set-TextForced "D22" '5.910'
Set-PlainText "E22" '  -1.63%  '

Set-PlainText "D23" '28.646.89'
Set-PlainText "E23" '  +0.94%  '

Set-TextForced "D24" '11.11'

Set-TextForced "D25" '2.096'
Set-PlainText "E25" '  -0.30%  '

Set-TextForced "D26" '159.00'
Set-PlainText "E26" '  +0.84%  '

Set-TextForced "D27" '20.46'
Set-PlainText "E27" '  -1.95%  '

Set-PlainText "D28" '2.007.37'
Set-PlainText "E28" '  -0.80%  '

Set-TextForced "D29" '2.344'
Set-PlainText "E29" '  -3.77%  '

Set-TextForced "D30" '122.68'
Set-PlainText "E30" '  -0.91%  '

Set-TextForced "D31" '1.105'
Set-PlainText "E31" '  -6.08%  '

Set-TextForced "D32" '0.1061'
Set-PlainText "E32" '  +0.60%  '

Set-TextForced "D33" '5.648'
Set-PlainText "E33" '  -2.40%  '

Set-TextForced "D34" '3.677'
Set-PlainText "E34" '  +1.17%  '

Set-TextForced "D35" '0.2243'
Set-PlainText "E35" '  +3.65%  '

Set-TextForced "D36" '0.06388'
Set-PlainText "E36" '  +4.82%  '

Set-TextForced "D37" '0.02300'
Set-PlainText "E37" '  -0.72%  '

Set-TextForced "D38" '8.769'
Set-PlainText "E38" '  -2.01%  '

Set-TextForced "D39" '5.034'
Set-PlainText "E39" '  -0.22%  '

Set-PlainText "B40" 'TrustWalletToken'
Set-PlainText "C40" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextForced "D40" '1.206'
Set-PlainText "E40" '  +4.48%  '

Set-PlainText "B41" 'Aptos'
Set-PlainText "C41" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextForced "D41" '11.25'
Set-PlainText "E41" '  -3.93%  '

Set-TextForced "D42" '0.6213'
Set-PlainText "E42" '  -3.24%  '

Set-TextForced "D43" '1.419'
Set-PlainText "E43" '  +2.81%  '

Set-TextForced "D44" '1.001'
Set-PlainText "E44" '  +0.02%  '

Set-TextForced "D45" '13.30'
Set-PlainText "E45" '  -0.61%  '

Set-PlainText "E46" '  -0.24%  '

Set-TextForced "D47" '0.5849'
Set-PlainText "E47" '  -2.42%  '

Set-TextForced "D48" '126.45'
Set-PlainText "E48" '  +3.66%  '

Set-TextForced "D49" '1.943'
Set-PlainText "E49" '  -0.13%  '

Set-TextForced "D50" '1.154'
Set-PlainText "E50" '  +0.49%  '

Set-TextForced "D51" '0.06890'
Set-PlainText "E51" '  +0.67%  '
